{"js": "// Unify/rewrite the CapacIT project description cell: collapse the\n// multiple draft paragraphs (original description, orange e-commerce\n// blurb, red draft notes, final orange summary, blank lines) into a\n// single paragraph with the new, unified wording.\n\nconst newText =\n  \"La plataforma web CapacIT ofrecer\u00e1 una amplia gama de cursos en l\u00ednea \" +\n  \"para la formaci\u00f3n en tecnolog\u00eda de la informaci\u00f3n (IT), desde cursos \" +\n  \"b\u00e1sicos para principiantes hasta cursos avanzados centrados en \u00e1reas \" +\n  \"de programaci\u00f3n y disciplinas relacionadas con IT. Estos cursos est\u00e1n \" +\n  \"dise\u00f1ados para satisfacer las necesidades tanto de aquellos que buscan \" +\n  \"adquirir habilidades t\u00e9cnicas para su carrera profesional, como de \" +\n  \"aquellos que buscan aprender por hobby o inter\u00e9s personal. La \" +\n  \"plataforma ofrecer\u00e1 cursos tanto gratuitos como de pago. Para acceder \" +\n  \"a los cursos de pago, se est\u00e1 desarrollando una pasarela de pago que \" +\n  \"permitir\u00e1 al usuario seleccionar los cursos que desea tomar y realizar \" +\n  \"el pago correspondiente de manera sencilla y segura, directamente en \" +\n  \"la p\u00e1gina web.\";\n\n// Locate the table/cell that holds the project-description paragraphs\n// (searching by content keeps this resilient to table-index drift).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet targetCell = null;\nfor (let i = 0; i < tables.items.length && !targetCell; i++) {\n  const table = tables.items[i];\n  table.load(\"values\");\n  await context.sync();\n  for (let r = 0; r < table.values.length && !targetCell; r++) {\n    for (let c = 0; c < table.values[r].length; c++) {\n      const cellText = table.values[r][c];\n      if (cellText && cellText.indexOf(\"CapacIT es una plataforma\") !== -1) {\n        targetCell = table.getCell(r, c);\n        break;\n      }\n    }\n  }\n}\n\nif (!targetCell) {\n  throw new Error(\"Could not locate the CapacIT description cell.\");\n}\n\nconst paragraphs = targetCell.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Replace the first paragraph's text in place (keeps its pPr/rPr, i.e.\n// no color override, sz/szCs = 24) with the unified text.\nparagraphs.items[0].insertText(newText, Word.InsertLocation.replace);\n\n// Remove every other paragraph in the cell (the old draft paragraphs).\nfor (let i = paragraphs.items.length - 1; i >= 1; i--) {\n  paragraphs.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Unify/rewrite the CapacIT project description cell: collapse the\n# multiple draft paragraphs (original description, orange e-commerce\n# blurb, red draft notes, final orange summary, blank lines) into a\n# single paragraph with the new, unified wording.\n\n$d = $word.ActiveDocument\n\n$newText = \"La plataforma web CapacIT ofrecer\u00e1 una amplia gama de cursos en l\u00ednea para la formaci\u00f3n en tecnolog\u00eda de la informaci\u00f3n (IT), desde cursos b\u00e1sicos para principiantes hasta cursos avanzados centrados en \u00e1reas de programaci\u00f3n y disciplinas relacionadas con IT. Estos cursos est\u00e1n dise\u00f1ados para satisfacer las necesidades tanto de aquellos que buscan adquirir habilidades t\u00e9cnicas para su carrera profesional, como de aquellos que buscan aprender por hobby o inter\u00e9s personal. La plataforma ofrecer\u00e1 cursos tanto gratuitos como de pago. Para acceder a los cursos de pago, se est\u00e1 desarrollando una pasarela de pago que permitir\u00e1 al usuario seleccionar los cursos que desea tomar y realizar el pago correspondiente de manera sencilla y segura, directamente en la p\u00e1gina web.\"\n\n$count = $d.Paragraphs.Count\n\n# Locate the first paragraph of the block (the long original description).\n$startIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text\n    if ($txt -like \"*CapacIT es una plataforma de formaci\u00f3n*\") {\n        $startIdx = $i\n        break\n    }\n}\nif ($startIdx -eq -1) {\n    throw \"Could not locate the start of the CapacIT description block.\"\n}\n\n# Locate the paragraph holding the last draft's closing sentence, then\n# include the two blank paragraphs that follow it (matches the original\n# cell layout: ...pago]Para/Web.[blank][blank]).\n$markerIdx = -1\nfor ($i = $startIdx; $i -le $count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text\n    if ($txt -like \"*directamente en la p\u00e1gina web*\") {\n        $markerIdx = $i\n        break\n    }\n}\nif ($markerIdx -eq -1) {\n    throw \"Could not locate the end marker of the CapacIT description block.\"\n}\n$endIdx = $markerIdx + 2\n\n# Replace the first paragraph's text in place (keeps its pPr/rPr, i.e. no\n# color override, sz/szCs = 24) with the unified text.\n$d.Paragraphs.Item($startIdx).Range.Text = $newText\n\n# Remove every other paragraph that belonged to the block. Deleting the\n# paragraph right after the (now-updated) first one repeatedly shifts\n# later ones down, so just keep removing the same index.\n$removeCount = $endIdx - $startIdx\nfor ($i = 1; $i -le $removeCount; $i++) {\n    $d.Paragraphs.Item($startIdx + 1).Range.Delete()\n}\n\nWrite-Output \"done\"\n"}
